$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped "cryptos" price table in place. Values are written as
# literal text (matching the workbook's existing inlineStr/shared-string
# cells) - price/percent columns are display strings, not live numbers.
# Cells whose new text happens to look like a plain number (e.g. "1.00",
# "0.0000171") are written with a leading apostrophe so Excel keeps storing
# them as text instead of silently re-typing them as numeric values.

# Row 2
$ws.Range('D2').Value = '60.665.79'
$ws.Range('E2').Value = '  -1.35%  '

# Row 3
$ws.Range('D3').Value = '3.386.38'
$ws.Range('E3').Value = '  -1.56%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').Value = '''568.28'
$ws.Range('E5').Value = '  -1.85%  '

# Row 6
$ws.Range('D6').Value = '''141.25'
$ws.Range('E6').Value = '  -4.54%  '

# Row 7
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('D8').Value = '3.388.17'
$ws.Range('E8').Value = '  -1.54%  '

# Row 9
$ws.Range('D9').Value = '''0.473'
$ws.Range('E9').Value = '  -0.35%  '

# Row 10
$ws.Range('D10').Value = '''7.52'
$ws.Range('E10').Value = '  -2.78%  '

# Row 11
$ws.Range('D11').Value = '''0.124'
$ws.Range('E11').Value = '  -1.61%  '

# Row 12
$ws.Range('D12').Value = '''0.388'
$ws.Range('E12').Value = '  -0.35%  '

# Row 13
$ws.Range('D13').Value = '3.965.23'
$ws.Range('E13').Value = '  -1.60%  '

# Row 14
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''28.27'
$ws.Range('E14').Value = '  +1.63%  '

# Row 15
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').Value = '''0.123'
$ws.Range('E15').Value = '  +0.96%  '

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000171'
$ws.Range('E16').Value = '  -1.54%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.387.31'
$ws.Range('E17').Value = '  -1.54%  '

# Row 18
$ws.Range('D18').Value = '60.803.31'
$ws.Range('E18').Value = '  -1.32%  '

# Row 19
$ws.Range('D19').Value = '''6.22'
$ws.Range('E19').Value = '  -0.61%  '

# Row 20
$ws.Range('D20').Value = '''13.88'
$ws.Range('E20').Value = '  -2.71%  '

# Row 21
$ws.Range('D21').Value = '''9.01'
$ws.Range('E21').Value = '  -4.80%  '

# Row 22
$ws.Range('D22').Value = '''388.34'
$ws.Range('E22').Value = '  +0.49%  '

# Row 23
$ws.Range('D23').Value = '''0.556'
$ws.Range('E23').Value = '  -1.43%  '

# Row 24
$ws.Range('D24').Value = '''73.16'
$ws.Range('E24').Value = '  +0.32%  '

# Row 25
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('D26').Value = '''0.0000117'
$ws.Range('E26').Value = '  -6.09%  '

# Row 27
$ws.Range('D27').Value = '3.520.57'
$ws.Range('E27').Value = '  -1.89%  '

# Row 28
$ws.Range('D28').Value = '''0.178'
$ws.Range('E28').Value = '  -1.21%  '

# Row 29
$ws.Range('E29').Value = '  +0.49%  '

# Row 30
$ws.Range('D30').Value = '''7.40'
$ws.Range('E30').Value = '  -3.53%  '

# Row 31
$ws.Range('E31').Value = '  -1.32%  '

# Row 32
$ws.Range('D32').Value = '''7.91'
$ws.Range('E32').Value = '  -3.52%  '

# Row 33
$ws.Range('D33').Value = '''1.40'
$ws.Range('E33').Value = '  -7.99%  '

# Row 34
$ws.Range('E34').Value = '  -0.04%  '

# Row 35
$ws.Range('D35').Value = '''23.46'
$ws.Range('E35').Value = '  -2.21%  '

# Row 36
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '''6.87'
$ws.Range('E36').Value = '  -1.43%  '

# Row 37
$ws.Range('D37').Value = '''167.80'
$ws.Range('E37').Value = '  +0.80%  '

# Row 38
$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').Value = '3.417.64'
$ws.Range('E38').Value = '  -1.55%  '

# Row 39
$ws.Range('E39').Value = '  -3.56%  '

# Row 40
$ws.Range('E40').Value = '  -4.02%  '

# Row 41
$ws.Range('D41').Value = '''0.0774'
$ws.Range('E41').Value = '  -0.82%  '

# Row 42
$ws.Range('D42').Value = '''27.35'
$ws.Range('E42').Value = '  +1.75%  '

# Row 43
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.781'
$ws.Range('E43').Value = '  -1.95%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  -0.03%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''41.55'
$ws.Range('E45').Value = '  -1.88%  '

# Row 46
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '''1.69'
$ws.Range('E46').Value = '  -0.60%  '

# Row 47
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '''4.39'
$ws.Range('E47').Value = '  -1.72%  '

# Row 48
$ws.Range('D48').Value = '2.546.49'
$ws.Range('E48').Value = '  -1.25%  '

# Row 49
$ws.Range('E49').Value = '  -2.61%  '

# Row 50
$ws.Range('D50').Value = '''23.00'
$ws.Range('E50').Value = '  -0.54%  '

# Row 51
$ws.Range('D51').Value = '''6.75'
$ws.Range('E51').Value = '  -1.89%  '
